$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.463.00'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -3.67%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.994.54'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -3.15%  '

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.15%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '536.56'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.61%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.01'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.92%  '

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '2.987.68'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.10%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.38%  '

# Row 10
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -5.76%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.10'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.82%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.445'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.27%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000221'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -3.20%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.66'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.71%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.490.91'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.76%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '61.603.05'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -3.38%  '

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.61%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.007.09'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -2.71%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.60'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.83%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '466.35'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -4.75%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.18'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.52%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.670'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -4.84%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.89'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.40%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.47'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.61%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.91'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.90%  '

# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.22%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.67'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.18%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.69'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -7.39%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.13%  '

# Row 30
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.98%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.52'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -2.94%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.87'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.06%  '

# Row 33
$ws.Range('B33').Value = 'Stacks'
$ws.Range('C33').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.27'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -6.08%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '55.23'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.22%  '

# Row 35
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.44'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.09%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.87'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.49%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '451.10'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -9.15%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.159.74'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.82%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0782'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.37%  '

# Row 40
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0383'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.20%  '

# Row 41
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.119'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.45%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.07'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.33%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.43'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -9.62%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '26.34'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.62%  '

# Row 45
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.10%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.242'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -7.06%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.98'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -5.31%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '118.22'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -3.31%  '

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.05%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0492'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -9.35%  '

# Row 51
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +6.33%  '
